$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-CellText($tbl, $row, $col, $newText) {
    $cell = $tbl.Cell($row, $col)
    $full = $cell.Range.Text
    $len = $full.Length - 2   # drop trailing paragraph mark + end-of-cell mark
    $start = $cell.Range.Start
    $rng = $d.Range($start, $start + $len)
    $rng.Text = $newText
}

function Get-CellText($tbl, $row, $col) {
    $cell = $tbl.Cell($row, $col)
    $full = $cell.Range.Text
    $len = $full.Length - 2
    return $full.Substring(0, $len)
}

# Rename header cells (row 1): GEAR.DIV.* columns lose the species qualifier,
# and the CATCH.DEP.* columns are relabeled to match the species each one
# actually reports on (cod vs hake).
Set-CellText $table 1 2 "GEAR.DIV.cod"
Set-CellText $table 1 3 "GEAR.DIV.hake"
Set-CellText $table 1 5 "CATCH.DEP.cod"
Set-CellText $table 1 6 "CATCH.DEP.hake"

# The CATCH.DEP.European.hake / CATCH.DEP.Atlantic.cod columns (5 and 6) had
# their values swapped between the two species - fix every data row.
$rowCount = $table.Rows.Count
for ($r = 2; $r -le $rowCount; $r++) {
    $left = Get-CellText $table $r 5
    $right = Get-CellText $table $r 6
    Set-CellText $table $r 5 $right
    Set-CellText $table $r 6 $left
}
